$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(3408, 45565, 36486.84, 24.88, 3.93, 1.15),
    @(3409, 45566, 36480.21, 24.87, 3.93, 1.15),
    @(3410, 45568, 35693.21, 24.34, 3.84, 1.18),
    @(3411, 45569, 35354.99, 24.11, 3.81, 1.19),
    @(3412, 45572, 34939.32, 23.82, 3.76, 1.21),
    @(3413, 45573, 35356.84, 24.11, 3.81, 1.19),
    @(3414, 45574, 35382.94, 24.12, 3.81, 1.19),
    @(3415, 45575, 35386.82, 24.13, 3.81, 1.19),
    @(3416, 45576, 35352.41, 24.1, 3.81, 1.19),
    @(3417, 45579, 35547.57, 24.24, 3.83, 1.19),
    @(3418, 45580, 35481.33, 24.18, 3.82, 1.19),
    @(3419, 45581, 35357.84, 24.1, 3.81, 1.19),
    @(3420, 45582, 34960.64, 23.83, 3.76, 1.2),
    @(3421, 45583, 35083.28, 23.93, 3.78, 1.21),
    @(3422, 45586, 34928.13, 23.75, 3.76, 1.2),
    @(3423, 45587, 34400.21, 23.3, 3.7, 1.23),
    @(3424, 45588, 34362.99, 23.27, 3.7, 1.23),
    @(3425, 45589, 34301.75, 23.21, 3.69, 1.23),
    @(3426, 45590, 33958.49, 22.97, 3.65, 1.19),
    @(3427, 45593, 34166.46, 23.14, 3.68, 1.19),
    @(3428, 45594, 34384.01, 23.41, 3.7, 1.23),
    @(3429, 45595, 34212.58, 23.3, 3.68, 1.24),
    @(3430, 45596, 34034.49, 23.14, 3.66, 1.27),
    @(3431, 45597, 34204.39, 23.26, 3.68, 1.26),
    @(3432, 45600, 33760.47, 22.93, 3.63, 1.24),
    @(3433, 45601, 34036.2, 23.1, 3.66, 1.26),
    @(3434, 45602, 34507.77, 23.42, 3.71, 1.26),
    @(3435, 45603, 34112.08, 23.15, 3.67, 1.28),
    @(3436, 45604, 33967.44, 22.87, 3.65, 1.29),
    @(3437, 45607, 33981.48, 22.86, 3.65, 1.27),
    @(3438, 45608, 33577.34, 22.41, 3.61, 1.29),
    @(3439, 45609, 33076.62, 22.12, 3.56, 1.3),
    @(3440, 45610, 33108.05, 22.13, 3.56, 1.32),
    @(3441, 45614, 33019.66, 22.06, 3.55, 1.31),
    @(3442, 45615, 33120.54, 22.14, 3.56, 1.32),
    @(3443, 45617, 32794.98, 21.92, 3.52, 1.29),
    @(3444, 45618, 33501.86, 22.4, 3.6, 1.26),
    @(3445, 45621, 33962.87, 22.7, 3.65, 1.25),
    @(3446, 45622, 33948.51, 22.69, 3.65, 1.25),
    @(3447, 45623, 34127.93, 22.81, 3.66, 1.24),
    @(3448, 45624, 33749, 22.56, 3.61, 1.25),
    @(3449, 45625, 34044.81, 22.76, 3.64, 1.24),
    @(3450, 45628, 34240.75, 22.89, 3.66, 1.24),
    @(3451, 45629, 34505.83, 23.06, 3.69, 1.23),
    @(3452, 45630, 34553.64, 23.1, 3.7, 1.22),
    @(3453, 45631, 34865.41, 23.31, 3.73, 1.21),
    @(3454, 45632, 34873.66, 23.31, 3.73, 1.21),
    @(3455, 45635, 34778.95, 23.25, 3.72, 1.22),
    @(3456, 45636, 34787.08, 23.25, 3.72, 1.22),
    @(3457, 45637, 34834.43, 23.28, 3.73, 1.21),
    @(3458, 45638, 34713.96, 23.2, 3.72, 1.2),
    @(3459, 45639, 34956.09, 23.37, 3.74, 1.2),
    @(3460, 45642, 34858.65, 23.3, 3.73, 1.2),
    @(3461, 45643, 34408.57, 23, 3.68, 1.21),
    @(3462, 45644, 34158.48, 22.83, 3.66, 1.22),
    @(3463, 45645, 33815.14, 22.6, 3.62, 1.24),
    @(3464, 45646, 33222.58, 22.21, 3.56, 1.26),
    @(3465, 45649, 33427.62, 22.34, 3.58, 1.25),
    @(3466, 45650, 33394.26, 22.32, 3.57, 1.26),
    @(3467, 45652, 33447.03, 22.35, 3.58, 1.26),
    @(3468, 45653, 33475.1, 22.37, 3.58, 1.25),
    @(3469, 45656, 33296.75, 22.25, 3.56, 1.25),
    @(3470, 45657, 33289.71, 22.31, 3.57, 1.25),
    @(3471, 45658, 33425.78, 22.4, 3.59, 1.24),
    @(3472, 45659, 34014.41, 22.79, 3.65, 1.22),
    @(3473, 45660, 33811.97, 22.66, 3.63, 1.23),
    @(3474, 45663, 33166.1, 22.22, 3.56, 1.25),
    @(3475, 45664, 33265.99, 22.29, 3.57, 1.25),
    @(3476, 45665, 33195.19, 22.24, 3.56, 1.25),
    @(3477, 45666, 32931.9, 22.07, 3.53, 1.26),
    @(3478, 45667, 32704.2, 21.91, 3.51, 1.27),
    @(3479, 45670, 32054.66, 21.47, 3.44, 1.29),
    @(3480, 45671, 32304.45, 21.63, 3.47, 1.28)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

Write-Output "Added $($rows.Count) rows (3408-3480)"
